$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9613333333333334
$ws.Range("C2").Value = 0.7873333333333333

$ws.Range("B3").Value = 0.9593333333333334
$ws.Range("C3").Value = 0.8113333333333334

$ws.Range("B4").Value = 0.956
$ws.Range("C4").Value = 0.7646666666666667

$ws.Range("B5").Value = 0.9633333333333334
$ws.Range("C5").Value = 0.802

$ws.Range("B6").Value = 0.9613333333333334
$ws.Range("C6").Value = 0.7733333333333333
